$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "Mele et al 2019"
$ws.Range("B15").Value = "optimal monpol under dgain learning isnt PLT, its inflation targeting"
$ws.Range("B15").WrapText = $true

$ws.Range("B16").Select()
